$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-13 10:20:39"
$zhcn.Range("E5").Value = "2016-03-13 10:20:39"
$zhcn.Range("H4").Value = "2016-03-13 10:20:59"
$zhcn.Range("H5").Value = "2016-03-13 10:20:59"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-13 10:20:45"
$dede.Range("E5").Value = "2016-03-13 10:20:45"
$dede.Range("H4").Value = "2016-03-13 10:21:07"
$dede.Range("H5").Value = "2016-03-13 10:21:07"
